$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47: Problem 51 - N-Queens (Hard / Not started -> "Bad" red) ---
$ws.Range("A47").Value = 51
$ws.Range("B47").Value = "N-Queens"
$ws.Range("C5").Copy()
$ws.Range("C47").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D4").Copy()
$ws.Range("D47").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 48: Problem 12 - Integer to Roman (Medium / "Neutral" yellow) ---
$ws.Range("A48").Value = 12
$ws.Range("B48").Value = "Integer to Roman"
$ws.Range("C3").Copy()
$ws.Range("C48").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D48").PasteSpecial(-4122)

# --- Row 49: Problem 11 - Container With Most Water (Medium / "Neutral" yellow) ---
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Container With Most Water"
$ws.Range("C3").Copy()
$ws.Range("C49").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D49").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Update the view state to match where the user left off scrolled/selected ---
$ws.Range("F47").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
